$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

# New match-day rows (fecha serial 45850 = 2025-07-12), appended after the
# existing data which ends at row 331.
$rows = @(
    @{ fecha=45850; jugador="Fabian Caicedo";           equipo="Amarillo"; posicion="Arquero";       goles=0; autogoles=0; arquero=$true;  goles_recibidos=3; amarillas=0; rojas=0; asistencias=0; penales=0 },
    @{ fecha=45850; jugador="Cesar Augusto Estrada";     equipo="Amarillo"; posicion="Delantero";     goles=1; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=0; rojas=0; asistencias=0; penales=0 },
    @{ fecha=45850; jugador="Andres Jurado";             equipo="Amarillo"; posicion="Delantero";     goles=2; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=0; rojas=0; asistencias=0; penales=0 },
    @{ fecha=45850; jugador="Andres Tangarife";          equipo="Amarillo"; posicion="Delantero";     goles=0; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=0; rojas=0; asistencias=1; penales=0 },
    @{ fecha=45850; jugador="Armando Murillo";           equipo="Amarillo"; posicion="Defensa";       goles=0; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=0; rojas=0; asistencias=1; penales=0 },
    @{ fecha=45850; jugador="Gember Marin Sarria";       equipo="Azul";     posicion="Arquero";       goles=0; autogoles=0; arquero=$true;  goles_recibidos=3; amarillas=0; rojas=0; asistencias=0; penales=0 },
    @{ fecha=45850; jugador="Alexander Uribe";           equipo="Azul";     posicion="Mediocampista"; goles=1; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=0; rojas=0; asistencias=2; penales=0 },
    @{ fecha=45850; jugador="David Fernando Velasco";    equipo="Azul";     posicion="Delantero";     goles=2; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=0; rojas=0; asistencias=0; penales=0 },
    @{ fecha=45850; jugador="Carlos Fernando Valencia";  equipo="Azul";     posicion="Delantero";     goles=0; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=0; rojas=0; asistencias=1; penales=0 },
    @{ fecha=45850; jugador="Esteban ";                  equipo="Azul";     posicion="Mediocampista"; goles=0; autogoles=0; arquero=$false; goles_recibidos=0; amarillas=1; rojas=0; asistencias=0; penales=0 }
)

$startRow = 332
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.fecha
    $ws.Cells.Item($r, 2).Value = $row.jugador
    $ws.Cells.Item($r, 3).Value = $row.equipo
    $ws.Cells.Item($r, 4).Value = $row.posicion
    $ws.Cells.Item($r, 5).Value = $row.goles
    $ws.Cells.Item($r, 6).Value = $row.autogoles
    $ws.Cells.Item($r, 7).Value = $row.arquero
    $ws.Cells.Item($r, 8).Value = $row.goles_recibidos
    $ws.Cells.Item($r, 9).Value = $row.amarillas
    $ws.Cells.Item($r, 10).Value = $row.rojas
    $ws.Cells.Item($r, 11).Value = $row.asistencias
    $ws.Cells.Item($r, 12).Value = $row.penales
}

# Mirror the author's final on-screen selection (rows below the new data)
# as captured in the saved workbook view state.
$ws.Range("A342:XFD408").Select()
